$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The survey respondent in row 7 originally typed "United States" for the
# Country column; correct it to "USA" to match the other respondents.
$ws.Range("B7").Value = "USA"
$ws.Range("B7").Font.Name = "Calibri"

# Leave the selection where the edit naturally lands after pressing Enter.
$ws.Range("B8").Select()
